$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 43 with the new "24/6/2025" report entry
$ws.Range("D43").Value = "24/6/2025"
$ws.Range("E43").Value = 297
$ws.Range("F43").Value = 629
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 1012
$ws.Range("J43").Value = "N/A"

# Update the view: scroll/selection as recorded when the author saved
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("L35").Select()
